$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Delete column E (IMEI)
$ws.Columns.Item(5).Delete()
# Now "token" is column H; insert a new column after H (before old I/J = longitude) for "mac"
$ws.Columns.Item(9).Insert()
# Delete the last column (autoSub), now column R
$ws.Columns.Item(18).Delete()
